$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'43.624.52"
$ws.Cells.Item(2, 5).Value = "'  +1.24%  "

$ws.Cells.Item(3, 4).Value = "'2.272.16"
$ws.Cells.Item(3, 5).Value = "'  +0.30%  "

$ws.Cells.Item(4, 5).Value = "'  -0.12%  "

$ws.Cells.Item(5, 4).Value = "'119.01"
$ws.Cells.Item(5, 5).Value = "'  +8.09%  "

$ws.Cells.Item(6, 4).Value = "'267.05"
$ws.Cells.Item(6, 5).Value = "'  +1.05%  "

$ws.Cells.Item(7, 5).Value = "'  +4.18%  "

$ws.Cells.Item(8, 4).Value = "'1.00"
$ws.Cells.Item(8, 5).Value = "'  +0.15%  "

$ws.Cells.Item(9, 4).Value = "'0.620"
$ws.Cells.Item(9, 5).Value = "'  +3.24%  "

$ws.Cells.Item(10, 4).Value = "'47.47"
$ws.Cells.Item(10, 5).Value = "'  +0.61%  "

$ws.Cells.Item(11, 4).Value = "'0.0945"
$ws.Cells.Item(11, 5).Value = "'  +2.20%  "

$ws.Cells.Item(12, 4).Value = "'9.49"
$ws.Cells.Item(12, 5).Value = "'  +8.40%  "

$ws.Cells.Item(13, 4).Value = "'0.106"
$ws.Cells.Item(13, 5).Value = "'  -1.00%  "

$ws.Cells.Item(14, 4).Value = "'15.73"
$ws.Cells.Item(14, 5).Value = "'  +2.37%  "

$ws.Cells.Item(15, 4).Value = "'0.899"
$ws.Cells.Item(15, 5).Value = "'  +6.22%  "

$ws.Cells.Item(16, 4).Value = "'2.615.11"
$ws.Cells.Item(16, 5).Value = "'  +0.23%  "

$ws.Cells.Item(17, 4).Value = "'2.274.65"
$ws.Cells.Item(17, 5).Value = "'  +0.39%  "

$ws.Cells.Item(18, 4).Value = "'43.525.88"
$ws.Cells.Item(18, 5).Value = "'  +1.25%  "

$ws.Cells.Item(19, 5).Value = "'  +2.32%  "

$ws.Cells.Item(20, 4).Value = "'6.93"
$ws.Cells.Item(20, 5).Value = "'  +1.94%  "

$ws.Cells.Item(21, 5).Value = "'  +1.82%  "

$ws.Cells.Item(22, 5).Value = "'  -2.47%  "

$ws.Cells.Item(23, 4).Value = "'234.44"
$ws.Cells.Item(23, 5).Value = "'  +1.62%  "

$ws.Cells.Item(24, 4).Value = "'2.92"
$ws.Cells.Item(24, 5).Value = "'  +3.06%  "

$ws.Cells.Item(25, 4).Value = "'9.52"
$ws.Cells.Item(25, 5).Value = "'  -0.73%  "

$ws.Cells.Item(26, 4).Value = "'12.38"
$ws.Cells.Item(26, 5).Value = "'  +10.09%  "

$ws.Cells.Item(27, 5).Value = "'  +1.90%  "

$ws.Cells.Item(28, 4).Value = "'41.98"
$ws.Cells.Item(28, 5).Value = "'  +4.71%  "

$ws.Cells.Item(29, 4).Value = "'3.33"
$ws.Cells.Item(29, 5).Value = "'  +1.63%  "

$ws.Cells.Item(30, 5).Value = "'  +0.25%  "

$ws.Cells.Item(31, 4).Value = "'174.43"
$ws.Cells.Item(31, 5).Value = "'  +1.70%  "

$ws.Cells.Item(32, 4).Value = "'21.44"
$ws.Cells.Item(32, 5).Value = "'  +1.10%  "

$ws.Cells.Item(33, 4).Value = "'0.0916"
$ws.Cells.Item(33, 5).Value = "'  +1.62%  "

$ws.Cells.Item(34, 5).Value = "'  +0.48%  "

$ws.Cells.Item(35, 5).Value = "'  +2.95%  "

$ws.Cells.Item(36, 4).Value = "'4.27"
$ws.Cells.Item(36, 5).Value = "'  +11.82%  "

$ws.Cells.Item(37, 5).Value = "'  +9.39%  "

$ws.Cells.Item(38, 4).Value = "'4.60"
$ws.Cells.Item(38, 5).Value = "'  -0.72%  "

$ws.Cells.Item(39, 5).Value = "'  +4.31%  "

$ws.Cells.Item(40, 4).Value = "'2.56"
$ws.Cells.Item(40, 5).Value = "'  -2.39%  "

$ws.Cells.Item(41, 4).Value = "'13.90"
$ws.Cells.Item(41, 5).Value = "'  +0.52%  "

$ws.Cells.Item(42, 2).Value = "'Algorand"
$ws.Cells.Item(42, 3).Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(42, 4).Value = "'0.240"
$ws.Cells.Item(42, 5).Value = "'  +2.52%  "

$ws.Cells.Item(43, 2).Value = "'MultiversX"
$ws.Cells.Item(43, 3).Value = "'https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Cells.Item(43, 4).Value = "'72.51"
$ws.Cells.Item(43, 5).Value = "'  -3.33%  "

$ws.Cells.Item(44, 5).Value = "'  -0.01%  "

$ws.Cells.Item(45, 5).Value = "'  +1.77%  "

$ws.Cells.Item(46, 4).Value = "'5.75"
$ws.Cells.Item(46, 5).Value = "'  -5.23%  "

$ws.Cells.Item(47, 4).Value = "'75.08"
$ws.Cells.Item(47, 5).Value = "'  +42.11%  "

$ws.Cells.Item(48, 2).Value = "'TrustWalletToken"
$ws.Cells.Item(48, 3).Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(48, 4).Value = "'1.27"
$ws.Cells.Item(48, 5).Value = "'  +2.68%  "

$ws.Cells.Item(49, 2).Value = "'Aave"
$ws.Cells.Item(49, 3).Value = "'https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(49, 4).Value = "'103.13"
$ws.Cells.Item(49, 5).Value = "'  +2.72%  "

$ws.Cells.Item(50, 4).Value = "'8.57"
$ws.Cells.Item(50, 5).Value = "'  -0.18%  "

$ws.Cells.Item(51, 4).Value = "'0.658"
$ws.Cells.Item(51, 5).Value = "'  +17.89%  "

